$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps its label text (the source re-emitted an extra duplicate shared
# string for it upstream, but the visible content is unchanged).
$ws.Range("A1").Value = "HK_G_acc_LG"

# A2:A49 are rescaled to the new distribution's values.
$ws.Range("A2").Value = 71.73913043478261
$ws.Range("A3").Value = 71.73913043478261
$ws.Range("A4").Value = 71.73913043478261
$ws.Range("A5").Value = 72.01086956521739
$ws.Range("A6").Value = 72.01086956521739
$ws.Range("A7").Value = 72.01086956521739
$ws.Range("A8").Value = 72.01086956521739
$ws.Range("A9").Value = 72.01086956521739
$ws.Range("A10").Value = 72.01086956521739
$ws.Range("A11").Value = 72.01086956521739
$ws.Range("A12").Value = 71.73913043478261
$ws.Range("A13").Value = 72.01086956521739
$ws.Range("A14").Value = 72.01086956521739
$ws.Range("A15").Value = 72.01086956521739
$ws.Range("A16").Value = 72.01086956521739
$ws.Range("A17").Value = 72.01086956521739
$ws.Range("A18").Value = 72.01086956521739
$ws.Range("A19").Value = 72.01086956521739
$ws.Range("A20").Value = 72.01086956521739
$ws.Range("A21").Value = 72.01086956521739
$ws.Range("A22").Value = 71.73913043478261
$ws.Range("A23").Value = 72.01086956521739
$ws.Range("A24").Value = 71.73913043478261
$ws.Range("A25").Value = 71.73913043478261
$ws.Range("A26").Value = 73.36956521739131
$ws.Range("A27").Value = 72.01086956521739
$ws.Range("A28").Value = 73.36956521739131
$ws.Range("A29").Value = 71.73913043478261
$ws.Range("A30").Value = 71.46739130434783
$ws.Range("A31").Value = 71.73913043478261
$ws.Range("A32").Value = 71.73913043478261
$ws.Range("A33").Value = 72.01086956521739
$ws.Range("A34").Value = 72.01086956521739
$ws.Range("A35").Value = 72.82608695652173
$ws.Range("A36").Value = 71.46739130434783
$ws.Range("A37").Value = 71.73913043478261
$ws.Range("A38").Value = 71.46739130434783
$ws.Range("A39").Value = 73.09782608695652
$ws.Range("A40").Value = 73.36956521739131
$ws.Range("A41").Value = 72.01086956521739
$ws.Range("A42").Value = 72.01086956521739
$ws.Range("A43").Value = 71.73913043478261
$ws.Range("A44").Value = 72.01086956521739
$ws.Range("A45").Value = 72.01086956521739
$ws.Range("A46").Value = 72.01086956521739
$ws.Range("A47").Value = 72.01086956521739
$ws.Range("A48").Value = 72.01086956521739
$ws.Range("A49").Value = 72.01086956521739
